# Generate Report for Handback
# - Overview sheet: Status columns (zh-cn/de-de) move from "Ready for handoff"
#   to "Handed back: in sync with en-US".
# - Each language sheet (zh-cn, de-de) gets its "Status" column updated the
#   same way, plus two new populated columns: "Latest Target File" (F) and
#   "Latest Handback File" (G), each holding a hyperlink to the file that was
#   handed back (mirroring the existing Source File Name / Latest Handoff
#   File hyperlinks), and the "Latest Handback DateTime" (H) column is
#   stamped with the handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: refresh the Status columns for both rows -------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- Per-language detail sheets --------------------------------------------
$mdTarget2  = "https://github.com/OpenLocalizationTest/oltest/blob/0f03600537a7e3fd8c1439786d0a95320e931dc4/e2e/aa846359-e5b7-4b1a-992e-45eab8e66c07.md"
$mdDisplay  = "aa846359-e5b7-4b1a-992e-45eab8e66c07.md"

$languages = @(
    @{ Name = "zh-cn"; XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbfc5ba3efab2263e3c7a5f7d6c3b9b7b9b17972/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/aa846359-e5b7-4b1a-992e-45eab8e66c07.6048720f6a2b7d4f85484e10c32f8fbb6ee781b6.zh-cn.xlf"; XlfDisplay = "aa846359-e5b7-4b1a-992e-45eab8e66c07.6048720f6a2b7d4f85484e10c32f8fbb6ee781b6.zh-cn.xlf"; HandbackTime = "2016-03-11 22:44:43" },
    @{ Name = "de-de"; XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/39b73928b3fce29ea80c77e0a3e831d0d53da03f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/aa846359-e5b7-4b1a-992e-45eab8e66c07.6048720f6a2b7d4f85484e10c32f8fbb6ee781b6.de-de.xlf"; XlfDisplay = "aa846359-e5b7-4b1a-992e-45eab8e66c07.6048720f6a2b7d4f85484e10c32f8fbb6ee781b6.de-de.xlf"; HandbackTime = "2016-03-11 22:44:49" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Name)

    foreach ($row in @(2, 3)) {
        # Status column
        $ws.Range("C$row").Value = $newStatus

        # Latest Target File (F) - hyperlink to the handed-off source .md
        $ws.Range("F$row").Value = $mdDisplay
        $ws.Hyperlinks.Add($ws.Range("F$row"), $mdTarget2, "", "", $mdDisplay) | Out-Null

        # Latest Handback File (G) - hyperlink to the translated .xlf
        $ws.Range("G$row").Value = $lang.XlfDisplay
        $ws.Hyperlinks.Add($ws.Range("G$row"), $lang.XlfTarget, "", "", $lang.XlfDisplay) | Out-Null

        # Latest Handback DateTime (H)
        $ws.Range("H$row").Value = $lang.HandbackTime
    }
}
